$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Adeyy"
$ws.Range("A3").Value = "Poef"
$ws.Range("A4").Value = "Lofd"

$ws.Range("E2").Value = "crewr354"
$ws.Range("E3").Value = "ftry3"
$ws.Range("E4").Value = "suhn35"

$ws.Range("C2").Select()
